$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from an existing header cell (H1) so the new
# headers get the same bold/border/centered formatting.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for columns I (I0) and J (IF)
$data = @(
    @(15, 16),
    @(6, 7),
    @(2, 3),
    @(9, 9),
    @(4, 5),
    @(3, 5),
    @(8, 9),
    @(5, 6),
    @(6, 7),
    @(6, 6),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
